$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: clone the number-format/border/fill (but not the value) of a
# reference cell that already carries the desired cellXf, then stamp the
# new value in. This reproduces the exact style index (s="...") used in the
# target workbook instead of Excel fabricating a brand-new style record.
# ---------------------------------------------------------------------------
function Set-GradeCell {
    param(
        [string]$Target,
        [string]$StyleSource,
        $Value
    )
    if ($StyleSource -ne $null) {
        $ws.Range($StyleSource).Copy()
        $ws.Range($Target).PasteSpecial(-4122)  # xlPasteFormats
    }
    $ws.Range($Target).Value = $Value
}

$excel.CutCopyMode = 0

# --- Row 10 (Бернакевич Елена): fill in all 7 homework grades -------------
Set-GradeCell "C10" "C5"  5
Set-GradeCell "D10" "C5"  5
Set-GradeCell "E10" "C5"  5
Set-GradeCell "F10" "C5"  5
Set-GradeCell "G10" "C5"  5
Set-GradeCell "H10" "C5"  5
Set-GradeCell "I10" "C5"  5

# --- Row 13 (Дубровская Мария): fill in all 7 homework grades -------------
Set-GradeCell "C13" "C5"  5
Set-GradeCell "D13" "C5"  5
Set-GradeCell "E13" "C5"  5
Set-GradeCell "F13" "C5"  5
Set-GradeCell "G13" "C5"  5
Set-GradeCell "H13" "C5"  5
Set-GradeCell "I13" "C5"  5

# --- Row 15 (Каиров Давид): replace the stray whitespace text with grades -
Set-GradeCell "C15" "C5"  5
Set-GradeCell "D15" "C5"  5
Set-GradeCell "E15" "C5"  5
Set-GradeCell "F15" "G25" 5
Set-GradeCell "G15" "I27" 5
Set-GradeCell "H15" "I27" 5

# --- Row 17 (Муллаянова Карина): extend grades into F:I --------------------
Set-GradeCell "F17" "F6" 5
Set-GradeCell "G17" "F6" 5
Set-GradeCell "H17" "F6" 5
Set-GradeCell "I17" "F6" 5

# --- Row 20 (Рогозин Даниил): add grades in G:H -----------------------------
Set-GradeCell "G20" $null 5
Set-GradeCell "H20" "I27" 5

# --- Row 29 (Султанов Денис): fill in all 7 homework grades ----------------
Set-GradeCell "C29" "C5"  5
Set-GradeCell "D29" "C5"  5
Set-GradeCell "E29" "C5"  5
Set-GradeCell "F29" "C5"  5
Set-GradeCell "G29" "C5"  5
Set-GradeCell "H29" "C5"  5
Set-GradeCell "I29" "C5"  5

# --- Update the sheet's saved selection/active cell -------------------------
$ws.Range("C29:I29").Select()

$wb.Application.CalculateFull()
